# Backlog.xlsx update:
#  - Add Kommentar (column G) notes explaining several backlog rows.
#  - Swap the Uppgift/Typ/Prioritering/Status (B:E) content of rows 29 and 30
#    ("Php sessioner" <-> "Admin funktionalitet"), carrying the Status-column
#    color formatting along with the data.
#  - Widen column B and update the view's selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- New comments in column G --------------------------------------------
# (Written in this order so the shared-string table is appended in the same
# sequence as the authored workbook.)
$ws.Range("G34").Value = "Mer enhetlig design och"
$ws.Range("G24").Value = "Alla beställningar och information om en användare ska tas bort"
$ws.Range("G25").Value = "Möjligheten att välja antal av en vara som ska läggas i varukorgen"
$ws.Range("G27").Value = "En sida där alla produkter inom en kategori visas som en lista."
$ws.Range("G32").Value = "Felmeddelanden ska ges felaktig inloggning och om recenstion läggs på ej köpt vara."
$ws.Range("G33").Value = "En siffra som visar hur många varor man har i kundvagnen."

# --- Swap rows 29 and 30 (B:D values, E fill color), row 29 <-> row 30 ----
# Column B/C/D keep their existing cell style, so swap the *values* only
# (Value2 snapshots the current content; .Value returns a live accessor in
# this host, so it must not be used for a read-then-write round trip).
$bTmp = $ws.Range("B29").Value2
$cTmp = $ws.Range("C29").Value2
$dTmp = $ws.Range("D29").Value2

$ws.Range("B29").Value = $ws.Range("B30").Value2
$ws.Range("C29").Value = $ws.Range("C30").Value2
$ws.Range("D29").Value = $ws.Range("D30").Value2

$ws.Range("B30").Value = $bTmp
$ws.Range("C30").Value = $cTmp
$ws.Range("D30").Value = $dTmp

# Column E carries no value, only the Status fill/style, so swap that via a
# scratch cell outside the used range.
$ws.Range("E29").Copy($ws.Range("E40")) | Out-Null
$ws.Range("E30").Copy($ws.Range("E29")) | Out-Null
$ws.Range("E40").Copy($ws.Range("E30")) | Out-Null
$ws.Range("E40").Clear() | Out-Null

# --- View/column-width touch ups ------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 42.33
$ws.Range("B31").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
